# Updates the cryptos price list (Coin/Link/Price/Volume(1h)) to the latest
# scraped snapshot. Rows 38/39 also swap: Aave now outranks ImmutableX.
#
# Price values in column D are stored as plain text (the sheet mixes
# formats like "31.237.76", "0.0641" and "0.0₃0722"), so for any value
# that looks numeric we force text by prefixing with an apostrophe and
# then reset the style back to Normal - this avoids Excel's automatic
# "typed value -> number" conversion while leaving no stray number-format
# style behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='31.237.76'; E='  +4.50%  '},
    @{Row=3; D='1.703.06'; E='  +4.13%  '},
    @{Row=4; D='0.998'; E='  -0.20%  '},
    @{Row=5; D='221.38'; E='  +2.76%  '},
    @{Row=6; E='  +2.63%  '},
    @{Row=7; D='0.997'; E='  -0.24%  '},
    @{Row=8; D='29.82'; E='  +3.56%  '},
    @{Row=9; D='45.22'; E='  +3.17%  '},
    @{Row=10; E='  +3.27%  '},
    @{Row=11; D='0.0641'; E='  +5.11%  '},
    @{Row=12; E='  +1.18%  '},
    @{Row=13; D='1.942.06'; E='  +3.89%  '},
    @{Row=14; D='1.698.22'; E='  +3.86%  '},
    @{Row=15; D='10.35'; E='  +9.59%  '},
    @{Row=16; D='0.615'; E='  +5.07%  '},
    @{Row=17; D='4.14'; E='  +6.95%  '},
    @{Row=18; D='31.216.93'; E='  +4.37%  '},
    @{Row=19; D='67.12'; E='  +3.64%  '},
    @{Row=20; D='249.21'; E='  +3.77%  '},
    @{Row=21; D='0.0₃0722'; E='  +2.70%  '},
    @{Row=22; D='0.999'; E='  -0.10%  '},
    @{Row=23; D='4.29'; E='  +3.76%  '},
    @{Row=24; D='10.15'; E='  +2.42%  '},
    @{Row=25; E='  -1.15%  '},
    @{Row=26; D='158.46'; E='  +0.66%  '},
    @{Row=27; D='15.99'; E='  +3.03%  '},
    @{Row=28; D='0.113'; E='  +3.13%  '},
    @{Row=29; D='6.75'; E='  +1.76%  '},
    @{Row=30; D='0.997'; E='  -0.26%  '},
    @{Row=31; D='0.0503'; E='  +2.96%  '},
    @{Row=32; D='3.65'; E='  +7.66%  '},
    @{Row=33; E='  +3.76%  '},
    @{Row=34; D='3.38'; E='  +5.83%  '},
    @{Row=35; D='1.524.06'; E='  +7.17%  '},
    @{Row=36; E='  +3.10%  '},
    @{Row=37; E='  +1.36%  '},
    @{Row=38; B='Aave'; C='https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D='83.60'; E='  +9.94%  '},
    @{Row=39; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.617'; E='  +10.31%  '},
    @{Row=40; E='  +4.78%  '},
    @{Row=41; D='2.70'; E='  -2.47%  '},
    @{Row=42; E='  +0.33%  '},
    @{Row=43; D='2.05'; E='  +3.10%  '},
    @{Row=44; D='0.852'; E='  +2.27%  '},
    @{Row=45; D='0.0505'; E='  +1.00%  '},
    @{Row=46; E='  +3.20%  '},
    @{Row=47; D='0.999'; E='  -0.15%  '},
    @{Row=48; D='52.30'; E='  +7.26%  '},
    @{Row=49; E='  +4.53%  '},
    @{Row=50; D='1.826.70'; E='  +2.79%  '},
    @{Row=51; D='94.25'; E='  +1.37%  '}
)

function Set-TextCell {
    param($Range, [string]$Value)

    $looksNumeric = $Value -match '^[+-]?[0-9]*\.?[0-9]+$'

    if ($looksNumeric) {
        # Forces Excel to treat the entry as literal text (quote-prefix),
        # matching the sheet's inline-string storage, then strips the
        # resulting cell style back to Normal so no stray number format
        # is left registered on the cell.
        $Range.Value = "'" + $Value
        $Range.Style = "Normal"
    }
    else {
        $Range.Value = $Value
    }
}

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$row").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$row").Value = $u.C }
    if ($u.ContainsKey('D')) { Set-TextCell $ws.Range("D$row") $u.D }
    if ($u.ContainsKey('E')) { $ws.Range("E$row").Value = $u.E }
}
